# Auto-generated edit script applying the Hades_Profits.xlsx diff
# to the appropriate worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Range("H112").Value = 1664.5122
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 1724.2307
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 5172.6921
$ws.Range("M112").Value = -392
$ws.Range("N112").Value = -7388.6921

# Row 132
$ws.Range("H132").Value = 2581302
$ws.Range("I132").Value = 2485.5557
$ws.Range("J132").Value = 49000000
$ws.Range("K132").Value = 7456.6671
$ws.Range("L132").Value = 147000000
$ws.Range("M132").Value = -4926.6671
$ws.Range("N132").Value = -147005060

# Row 137
$ws.Range("H137").Value = 3336029.8
$ws.Range("I137").Value = 7145300.5
$ws.Range("J137").Value = 2918
$ws.Range("K137").Value = 21435901.5
$ws.Range("L137").Value = 8754
$ws.Range("M137").Value = -21433351.5
$ws.Range("N137").Value = -13854

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7779052
$ws.Range("I32").Value = 9836692
$ws.Range("J32").Value = 5744.778
$ws.Range("K32").Value = 9836692
$ws.Range("L32").Value = 5744.778
$ws.Range("M32").Value = -9836405
$ws.Range("N32").Value = -6318.778

# Row 45
$ws.Range("H45").Value = 4522.04
$ws.Range("I45").Value = 4472.227
$ws.Range("J45").Value = 4887.3335
$ws.Range("K45").Value = 4472.227
$ws.Range("L45").Value = 4887.3335
$ws.Range("M45").Value = -4095.227
$ws.Range("N45").Value = -5641.3335

# Row 61
$ws.Range("H61").Value = 47715984
$ws.Range("I61").Value = 55612910
$ws.Range("J61").Value = 334433.34
$ws.Range("K61").Value = 55612910
$ws.Range("L61").Value = 334433.34
$ws.Range("M61").Value = -55612698
$ws.Range("N61").Value = -334857.34

# Row 74
$ws.Range("H74").Value = 16801058
$ws.Range("I74").Value = 27889690
$ws.Range("J74").Value = 168108.33
$ws.Range("K74").Value = 27889690
$ws.Range("L74").Value = 168108.33
$ws.Range("M74").Value = -27888816
$ws.Range("N74").Value = -169856.33

# Row 77
$ws.Range("H77").Value = 16801058
$ws.Range("I77").Value = 27889690
$ws.Range("J77").Value = 168108.33
$ws.Range("K77").Value = 139448450
$ws.Range("L77").Value = 840541.6499999999
$ws.Range("M77").Value = -139444082
$ws.Range("N77").Value = -849277.6499999999

# Row 122
$ws.Range("H122").Value = 1351.5555
$ws.Range("I122").Value = 1342.8235
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4028.4705
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -1578.4705
$ws.Range("N122").Value = -9400

# Row 131
$ws.Range("H131").Value = 54285.715
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 54285.715
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 54285.715
$ws.Range("N131").Value = -64365.715

# Row 132
$ws.Range("H132").Value = 75968.22
$ws.Range("I132").Value = 49062.668
$ws.Range("J132").Value = 170137.67
$ws.Range("K132").Value = 147188.004
$ws.Range("L132").Value = 510413.01
$ws.Range("M132").Value = -144658.004
$ws.Range("N132").Value = -515473.01

# Row 136
$ws.Range("H136").Value = 47715984
$ws.Range("I136").Value = 55612910
$ws.Range("J136").Value = 334433.34
$ws.Range("K136").Value = 166838730
$ws.Range("L136").Value = 1003300.02
$ws.Range("M136").Value = -166836180
$ws.Range("N136").Value = -1008400.02

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 5622.3213
$ws.Range("I134").Value = 4892.7085
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 14678.1255
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -12143.1255
$ws.Range("N134").Value = -35070

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1767.9524
$ws.Range("I31").Value = 1606.65
$ws.Range("J31").Value = 4994
$ws.Range("K31").Value = 1606.65
$ws.Range("L31").Value = 4994
$ws.Range("M31").Value = -1311.65
$ws.Range("N31").Value = -5584

# Row 34
$ws.Range("H34").Value = 1767.9524
$ws.Range("I34").Value = 1606.65
$ws.Range("J34").Value = 4994
$ws.Range("K34").Value = 1606.65
$ws.Range("L34").Value = 4994
$ws.Range("M34").Value = -1404.65
$ws.Range("N34").Value = -5398

# Row 58
$ws.Range("H58").Value = 76928780
$ws.Range("I58").Value = 76928780
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 76928780
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -76928577

# Row 59
$ws.Range("H59").Value = 20000
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 20000
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 20000
$ws.Range("N59").Value = -22290

# Row 132
$ws.Range("H132").Value = 173458
$ws.Range("I132").Value = 11375
$ws.Range("J132").Value = 254499.5
$ws.Range("K132").Value = 34125
$ws.Range("L132").Value = 763498.5
$ws.Range("M132").Value = -31595
$ws.Range("N132").Value = -768558.5

# Row 134
$ws.Range("H134").Value = 65018.59
$ws.Range("I134").Value = 1612.8889
$ws.Range("J134").Value = 136350
$ws.Range("K134").Value = 4838.6667
$ws.Range("L134").Value = 409050
$ws.Range("M134").Value = -2303.6667
$ws.Range("N134").Value = -414120

# Row 135
$ws.Range("H135").Value = 55450
$ws.Range("I135").Value = 56000
$ws.Range("J135").Value = 54900
$ws.Range("K135").Value = 56000
$ws.Range("L135").Value = 54900
$ws.Range("M135").Value = -50930
$ws.Range("N135").Value = -65040

# Row 136
$ws.Range("H136").Value = 76928780
$ws.Range("I136").Value = 76928780
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 230786340
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -230783790

# Row 138
$ws.Range("H138").Value = 49040
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 49040
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 49040
$ws.Range("N138").Value = -59320

$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 1187.92
$ws.Range("I122").Value = 516.5
$ws.Range("J122").Value = 1399.9474
$ws.Range("K122").Value = 4648.5
$ws.Range("L122").Value = 12599.5266
$ws.Range("M122").Value = -2198.5
$ws.Range("N122").Value = -17499.5266

# Row 123
$ws.Range("H123").Value = 3446.6
$ws.Range("I123").Value = 1000
$ws.Range("J123").Value = 4058.25
$ws.Range("K123").Value = 3000
$ws.Range("L123").Value = 12174.75
$ws.Range("M123").Value = -550
$ws.Range("N123").Value = -17074.75

# Row 132
$ws.Range("H132").Value = 2125
$ws.Range("I132").Value = 4243.5
$ws.Range("J132").Value = 1654.2222
$ws.Range("K132").Value = 38191.5
$ws.Range("L132").Value = 14887.9998
$ws.Range("M132").Value = -35661.5
$ws.Range("N132").Value = -19947.9998

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 2251.25
$ws.Range("I97").Value = 2251.25
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2251.25
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1755.25
$ws.Range("N97").ClearContents()

# Row 102
$ws.Range("H102").Value = 1663.6
$ws.Range("I102").Value = 1133.1111
$ws.Range("J102").Value = 2459.3333
$ws.Range("K102").Value = 1133.1111
$ws.Range("L102").Value = 2459.3333
$ws.Range("M102").Value = 488.8888999999999
$ws.Range("N102").Value = -5703.3333

# Row 107
$ws.Range("H107").Value = 2270.75
$ws.Range("I107").Value = 1690
$ws.Range("J107").Value = 2851.5
$ws.Range("K107").Value = 1690
$ws.Range("L107").Value = 2851.5
$ws.Range("M107").Value = 230
$ws.Range("N107").Value = -6691.5

# Row 132
$ws.Range("H132").Value = 57320.676
$ws.Range("I132").Value = 42907
$ws.Range("J132").Value = 87349.164
$ws.Range("K132").Value = 128721
$ws.Range("L132").Value = 262047.492
$ws.Range("M132").Value = -126191
$ws.Range("N132").Value = -267107.492

# Row 134
$ws.Range("H134").Value = 26663
$ws.Range("I134").Value = 29000
$ws.Range("J134").Value = 24326
$ws.Range("K134").Value = 87000
$ws.Range("L134").Value = 72978
$ws.Range("M134").Value = -84465
$ws.Range("N134").Value = -78048

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 8000
$ws.Range("I40").Value = 8000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 8000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -7864

# Row 135
$ws.Range("H135").Value = 32943
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 32943
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 32943
$ws.Range("N135").Value = -43083

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 56311.223
$ws.Range("I100").Value = 46082
$ws.Range("J100").Value = 72385.71000000001
$ws.Range("K100").Value = 92164
$ws.Range("L100").Value = 144771.42
$ws.Range("M100").Value = -91623
$ws.Range("N100").Value = -145853.42

# Row 113
$ws.Range("H113").Value = 807.5833
$ws.Range("I113").Value = 430.91666
$ws.Range("J113").Value = 1184.25
$ws.Range("K113").Value = 1292.74998
$ws.Range("L113").Value = 3552.75
$ws.Range("M113").Value = 877.2500199999999
$ws.Range("N113").Value = -7892.75
